$d = $word.ActiveDocument

# After the final inline picture's paragraph (anchorId 3B7E47D7 / docPr
# "Picture 4" -- the last InlineShape in the document), add a blank
# paragraph followed by a paragraph containing the repo URL, matching
# the end of the document body just before the section break.
$picPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$picPara.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter()

$linkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$linkPara.Range.Text = "https://github.com/MuriloManhas88/GsGovernanca"

# Word marks runs holding freshly inserted/refreshed pictures as
# NoProofing (<w:rPr><w:noProof/></w:rPr>), which this last picture was
# still missing. Set this after the paragraph inserts above so the new
# paragraphs don't inherit the NoProofing formatting.
$lastShape = $d.InlineShapes.Item($d.InlineShapes.Count)
$lastShape.Range.NoProofing = $true
